$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 holds the f6cd01df... handoff entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 16:52:41"
$wsZhCn.Range("H2").Value = "2016-03-17 16:52:59"

# de-de sheet: row 2 holds the f6cd01df... handoff entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 16:52:45"
$wsDeDe.Range("H2").Value = "2016-03-17 16:53:10"
